$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = "General"
}

function Set-NumCell($row, $col, $num) {
    $rng = $ws.Cells.Item($row, $col)
    $rng.Value = $num
}

function Clear-Cell($row, $col) {
    $ws.Cells.Item($row, $col).ClearContents()
}

# Row 8
Set-TextCell 8 1 ""
Set-NumCell 8 2 0
Set-NumCell 8 3 45
Set-TextCell 8 4 '2.0'
$ws.Cells.Item(8, 5).Value = 'Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it''s  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumCell 8 6 0
Set-TextCell 8 7 '0.00'
Set-NumCell 8 8 0
Set-TextCell 8 9 ""

# Row 9
$ws.Cells.Item(9, 1).Value = 'P. point'
Set-NumCell 9 2 0
Set-NumCell 9 3 84
Set-TextCell 9 4 '6'
$ws.Cells.Item(9, 5).Value = 'On board'
Set-NumCell 9 6 136
Set-TextCell 9 7 '11424.00'
Set-NumCell 9 8 0
Set-TextCell 9 9 ""

# Row 10
$ws.Cells.Item(10, 1).Value = 'Each'
Set-NumCell 10 2 0
Set-NumCell 10 3 57
Set-TextCell 10 4 '3.0'
$ws.Cells.Item(10, 5).Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumCell 10 6 23
Set-TextCell 10 7 '1311.00'
Set-NumCell 10 8 0
Set-TextCell 10 9 ""

# Row 11
$ws.Cells.Item(11, 1).Value = 'Each'
Set-NumCell 11 2 0
Set-NumCell 11 3 11
Set-TextCell 11 4 '4.0'
$ws.Cells.Item(11, 5).Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumCell 11 6 50
Set-TextCell 11 7 '550.00'
Set-NumCell 11 8 0
Set-TextCell 11 9 ""

# Row 12
$ws.Cells.Item(12, 1).Value = 'Each'
Set-NumCell 12 2 0
Set-NumCell 12 3 34
Set-TextCell 12 4 '6.0'
$ws.Cells.Item(12, 5).Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumCell 12 6 78
Set-TextCell 12 7 '2652.00'
Set-NumCell 12 8 0
Set-TextCell 12 9 ""

# Row 13
$ws.Cells.Item(13, 1).Value = 'Each'
Set-NumCell 13 2 0
Set-NumCell 13 3 51
Set-TextCell 13 4 '9.0'
$ws.Cells.Item(13, 5).Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumCell 13 6 219
Set-TextCell 13 7 '11169.00'
Set-NumCell 13 8 0
Set-TextCell 13 9 ""

# Row 14
$ws.Cells.Item(14, 1).Value = 'Each'
Set-NumCell 14 2 0
Set-NumCell 14 3 13
Set-TextCell 14 4 '10.0'
$ws.Cells.Item(14, 5).Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumCell 14 6 303
Set-TextCell 14 7 '3939.00'
Set-NumCell 14 8 0
Set-TextCell 14 9 ""

# Row 15
Set-TextCell 15 1 ""
Set-NumCell 15 2 0
Set-NumCell 15 3 8
Set-TextCell 15 4 '11.0'
$ws.Cells.Item(15, 5).Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumCell 15 6 0
Set-TextCell 15 7 '0.00'
Set-NumCell 15 8 0
Set-TextCell 15 9 ""

# Row 16
$ws.Cells.Item(16, 1).Value = 'R. mtr.'
Set-NumCell 16 2 0
Set-NumCell 16 3 57
Set-TextCell 16 4 '17'
$ws.Cells.Item(16, 5).Value = '25 mm'
Set-NumCell 16 6 56
Set-TextCell 16 7 '3192.00'
Set-NumCell 16 8 0
Set-TextCell 16 9 ""

# Row 17
$ws.Cells.Item(17, 1).Value = 'Mtr.'
Set-NumCell 17 2 0
Set-NumCell 17 3 74
Set-TextCell 17 4 '19'
$ws.Cells.Item(17, 5).Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
Set-NumCell 17 6 81
Set-TextCell 17 7 '5994.00'
Set-NumCell 17 8 0
Set-TextCell 17 9 ""

# Row 18
$ws.Cells.Item(18, 1).Value = 'Set'
Set-NumCell 18 2 0
Set-NumCell 18 3 98
Set-TextCell 18 4 '13.0'
$ws.Cells.Item(18, 5).Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
Set-NumCell 18 6 5733
Set-TextCell 18 7 '561834.00'
Set-NumCell 18 8 0
Set-TextCell 18 9 ""

# Row 19
Set-TextCell 19 1 ""
Set-NumCell 19 2 0
Set-NumCell 19 3 58
Set-TextCell 19 4 '14.0'
$ws.Cells.Item(19, 5).Value = 'Supply & Laying following size earth wire in horizontal or vertical run in ground/surface/recess including riveting, soldering, saddles,  making connection with GI/Cu purity purity >95%  thimble etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
Set-NumCell 19 6 0
Set-TextCell 19 7 '0.00'
Set-NumCell 19 8 0
Set-TextCell 19 9 ""

# Row 20
$ws.Cells.Item(20, 1).Value = 'Mtr.'
Set-NumCell 20 2 0
Set-NumCell 20 3 59
Set-TextCell 20 4 '23'
$ws.Cells.Item(20, 5).Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
Set-NumCell 20 6 20
Set-TextCell 20 7 '1180.00'
Set-NumCell 20 8 0
Set-TextCell 20 9 ""

# Row 21
Set-TextCell 21 1 ""
Set-NumCell 21 2 0
Set-NumCell 21 3 63
Set-TextCell 21 4 '15.0'
$ws.Cells.Item(21, 5).Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumCell 21 6 0
Set-TextCell 21 7 '0.00'
Set-NumCell 21 8 0
Set-TextCell 21 9 ""

# Row 22
$ws.Cells.Item(22, 1).Value = 'Each'
Set-NumCell 22 2 0
Set-NumCell 22 3 12
Set-TextCell 22 4 '25'
$ws.Cells.Item(22, 5).Value = '1200 mm Sweep BEE 1 Star rated (service value >=4.0 to < 4.5 )'
Set-NumCell 22 6 1890
Set-TextCell 22 7 '22680.00'
Set-NumCell 22 8 0
Set-TextCell 22 9 ""

# Row 23
$ws.Cells.Item(23, 1).Value = 'Each'
Set-NumCell 23 2 0
Set-NumCell 23 3 62
Set-TextCell 23 4 '27'
$ws.Cells.Item(23, 5).Value = '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
Set-NumCell 23 6 492
Set-TextCell 23 7 '30504.00'
Set-NumCell 23 8 0
Set-TextCell 23 9 ""

# Row 24
Set-TextCell 24 1 ""
Set-NumCell 24 2 0
Set-NumCell 24 3 29
Set-TextCell 24 4 '17.0'
$ws.Cells.Item(24, 5).Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumCell 24 6 0
Set-TextCell 24 7 '0.00'
Set-NumCell 24 8 0
Set-TextCell 24 9 ""

# Row 25
Set-TextCell 25 1 ""
Set-NumCell 25 2 0
Set-NumCell 25 3 31
Set-TextCell 25 4 '34'
$ws.Cells.Item(25, 5).Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
Set-NumCell 25 6 0
Set-TextCell 25 7 '0.00'
Set-NumCell 25 8 0
Set-TextCell 25 9 ""

# Row 26
$ws.Cells.Item(26, 1).Value = '%'
Set-NumCell 26 2 0
Set-NumCell 26 3 44
Set-TextCell 26 4 '37'
$ws.Cells.Item(26, 5).Value = 'Add Tender Premium '
Set-NumCell 26 6 0
Set-TextCell 26 7 '0.00'
Set-NumCell 26 8 0
Set-TextCell 26 9 ""

# Row 27
Set-TextCell 27 1 ""
Set-NumCell 27 2 0
Set-NumCell 27 3 4
Set-TextCell 27 4 '38'
$ws.Cells.Item(27, 5).Value = 'Grand Total'
Set-NumCell 27 6 0
Set-TextCell 27 7 '0.00'
Set-NumCell 27 8 0
Set-TextCell 27 9 ""

# Row 28
Set-TextCell 28 1 ""
Clear-Cell 28 2
Clear-Cell 28 3
Clear-Cell 28 4
Clear-Cell 28 5
Clear-Cell 28 6
Clear-Cell 28 7
Clear-Cell 28 8
Clear-Cell 28 9

# Row 29
Set-TextCell 29 1 ""
Set-TextCell 29 2 ""
Set-TextCell 29 3 ""
Set-TextCell 29 4 ""
$ws.Cells.Item(29, 5).Value = 'Grand Total Rs.'
Set-TextCell 29 6 ""
Set-TextCell 29 7 '656429.00'
Set-TextCell 29 8 '656429.00'
Set-TextCell 29 9 ""

# Row 30
Set-TextCell 30 1 ""
Set-TextCell 30 2 ""
Set-TextCell 30 3 ""
Set-TextCell 30 4 ""
$ws.Cells.Item(30, 5).Value = 'Tender Premium @ 0%'
Set-TextCell 30 6 ""
Set-TextCell 30 7 '0.00'
Set-TextCell 30 8 '0.00'
Set-TextCell 30 9 ""

# Row 31
Set-TextCell 31 1 ""
Set-TextCell 31 2 ""
Set-TextCell 31 3 ""
Set-TextCell 31 4 ""
$ws.Cells.Item(31, 5).Value = 'NET PAYABLE AMOUNT Rs.'
Set-TextCell 31 6 ""
Set-TextCell 31 7 '656429.00'
Set-TextCell 31 8 '656429.00'
Set-TextCell 31 9 ""

